$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.772.82"
$ws.Range("E2").Value = "  +1.53%  "
$ws.Range("D3").Value = "3.327.66"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "189.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "588.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("D12").Value = "3.904.93"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.03"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "68.861.85"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "3.327.29"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "447.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.98%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +2.52%  "
$ws.Range("D25").Value = "3.489.92"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("E26").Value = "  +2.51%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.190"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E29").Value = "  -0.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.58%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.17%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +5.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.62%  "
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.793"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("D43").Value = "2.694.99"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0682"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "330.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.57%  "
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.61%  "
$ws.Range("E51").Value = "  +3.42%  "
